# Daily automatic data refresh for epexspot_prices.xlsx
#  - "Prix Spot": a new daily column (01-dec) is inserted right before the
#    block of "-oct." columns (i.e. right where the new day's data belongs),
#    pushing the existing "01-oct." .. "31-oct." columns one slot to the
#    right. The brand-new day has no data yet, so every hourly row gets "-".
#  - "Gaz" / "CO2": two new trading days are appended at the bottom
#    (2025-11-29 and 2025-11-30).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": insert the new "01-dec" column at ED (column 134)
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Columns("ED:ED").Insert()

$wsPrix.Range("ED1").Value = "01-dec"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 134).Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the two new trading days
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A164:A165").NumberFormat = "@"
$wsGaz.Range("A164").Value = "2025-11-29"
$wsGaz.Range("B164").Value = 27.525
$wsGaz.Range("A165").Value = "2025-11-30"
$wsGaz.Range("B165").Value = 27.525
$wsGaz.Range("A164:A165").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "CO2": append the two new trading days (price not published yet)
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A164:A165").NumberFormat = "@"
$wsCO2.Range("A164").Value = "2025-11-29"
$wsCO2.Range("A165").Value = "2025-11-30"
$wsCO2.Range("A164:A165").Style = "Normal"
